# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp header
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 07:20"

# Row 4 - Estados Unidos: update stats
$ws.Range("B4").Value = 245341
$ws.Range("C4").Value = 464
$ws.Range("E4").Value = 228843
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 6095

# Row 23 - Australia: update stats
$ws.Range("E23").Value = 4701
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 28

# Row 36 - India: update stats
$ws.Range("B36").Value = 2567
$ws.Range("C36").Value = 24
$ws.Range("D36").Value = 192
$ws.Range("E36").Value = 2303

# Rows 39/40 - Tailandia moves above Arabia Saudita (Tailandia updated, Arabia Saudita unchanged but shifted down)
$ws.Range("A39").Value = "Tailandia"
$ws.Range("B39").Value = 1978
$ws.Range("C39").Value = 103
$ws.Range("D39").Value = 581
$ws.Range("E39").Value = 1378
$ws.Range("F39").Value = 23
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 19

$ws.Range("A40").Value = "Arabia Saudita"
$ws.Range("B40").Value = 1885
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 328
$ws.Range("E40").Value = 1536
$ws.Range("F40").Value = 31
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 21

# Row 46 - Sudafrica: update stats
$ws.Range("D46").Value = 95
$ws.Range("E46").Value = 1362

# Rows 67/68 - Lituania moves above Armenia (Lituania updated, Armenia unchanged but shifted down)
$ws.Range("A67").Value = "Lituania"
$ws.Range("B67").Value = 696
$ws.Range("C67").Value = 47
$ws.Range("D67").Value = 7
$ws.Range("E67").Value = 680
$ws.Range("F67").Value = 11
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 9

$ws.Range("A68").Value = "Armenia"
$ws.Range("B68").Value = 663
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 33
$ws.Range("E68").Value = 623
$ws.Range("F68").Value = 30
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 7

# Row 70 - Hungria: update stats
$ws.Range("B70").Value = 623
$ws.Range("C70").Value = 38
$ws.Range("D70").Value = 43
$ws.Range("E70").Value = 554
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 26

# Row 96 - Vietnam: update stats
$ws.Range("D96").Value = 85
$ws.Range("E96").Value = 148

# Row 99 - Uzbekistan: update stats
$ws.Range("B99").Value = 221
$ws.Range("C99").Value = 16
$ws.Range("E99").Value = 194

# Rows 116/117 - Kirguistan moves above Guadalupe (Kirguistan updated, Guadalupe unchanged but shifted down)
$ws.Range("A116").Value = "Kirguistan"
$ws.Range("B116").Value = 130
$ws.Range("C116").Value = 14
$ws.Range("D116").Value = 5
$ws.Range("E116").Value = 124
$ws.Range("F116").Value = 5
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 1

$ws.Range("A117").Value = "Guadalupe"
$ws.Range("B117").Value = 128
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 24
$ws.Range("E117").Value = 98
$ws.Range("F117").Value = 14
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 6
